$wb = $excel.ActiveWorkbook

$sheetNames = @("BPaFF-BITPTaP", "BPaFF-BDTPTPF")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Right-align the "Boolean" header cell (B1) - new style with horizontal="right"
    $ws.Range("B1").HorizontalAlignment = -4152  # xlRight

    # New row 15: crude oil -> references B11 (petroleum)
    $ws.Cells.Item(15, 1).Value = "crude oil"
    $ws.Cells.Item(15, 2).Formula = "=B11"

    # New row 16: heavy or residual fuel oil -> references B11 (petroleum)
    $ws.Cells.Item(16, 1).Value = "heavy or residual fuel oil"
    $ws.Cells.Item(16, 2).Formula = "=B11"

    # New row 17: municipal solid waste -> references B9 (biomass)
    $ws.Cells.Item(17, 1).Value = "municipal solid waste"
    $ws.Cells.Item(17, 2).Formula = "=B9"
}
